$d = $word.ActiveDocument

# Change 1: expand the sentence about where the key needs to be added
$d.Content.Find.Execute(
    "but the following key needs to be added to the Config/appSettings.config file:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "but the following key needs to be added to appSettings (in Web.config or Config/appSettings.config, depending on your Umbraco setup):",
    2
) | Out-Null

# Change 2: "...web site's Config/appSettings.config file). That means" -> "...web site's appSettings). That means"
$d.Content.Find.Execute(
    "can be added to the web site's Config/appSettings.config file). That means",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "can be added to the web site's appSettings). That means",
    2
) | Out-Null

# Change 3: "...by entering the following to the Config/appSettings.config file:" -> "...by entering the following to appSettings:"
$d.Content.Find.Execute(
    "by entering the following to the Config/appSettings.config file:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "by entering the following to appSettings:",
    2
) | Out-Null

# Change 4: "...web site's Config/appSettings.config file." -> "...web site's appSettings."
$d.Content.Find.Execute(
    "can be added to the web site's Config/appSettings.config file.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "can be added to the web site's appSettings.",
    2
) | Out-Null

# Change 5: "...web site's Config/appSettings.config file. Then, the languages" -> "...web site's appSettings. Then, the languages"
$d.Content.Find.Execute(
    "needs to be entered into the web site's Config/appSettings.config file. Then, the languages",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "needs to be entered into the web site's appSettings. Then, the languages",
    2
) | Out-Null
